$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.022.15"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.912.86"
$ws.Range("E3").Value = "  +0.07%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "589.41"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.42%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "144.92"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.16%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "6.89"
$c.ClearFormats()
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -0.13%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "33.45"
$c.ClearFormats()
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D15").Value = "3.395.94"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "60.902.33"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "2.911.06"
$ws.Range("E18").Value = "  -0.03%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "432.59"
$c.ClearFormats()
$ws.Range("E19").Value = "  +1.08%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.39"
$c.ClearFormats()
$ws.Range("E20").Value = "  -1.02%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.676"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.22%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "81.46"
$c.ClearFormats()
$ws.Range("E23").Value = "  +1.10%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.82"
$c.ClearFormats()
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  -1.21%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.78"
$c.ClearFormats()
$ws.Range("E26").Value = "  -1.18%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.ClearFormats()
$ws.Range("E28").Value = "  +5.04%  "
$ws.Range("E29").Value = "  -0.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.97"
$c.ClearFormats()
$ws.Range("E30").Value = "  -3.13%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "26.53"
$c.ClearFormats()
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "0.0₃0869"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("E35").Value = "  -0.18%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.62"
$c.ClearFormats()
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -3.07%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.55"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  -4.20%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "40.85"
$c.ClearFormats()
$ws.Range("E42").Value = "  -1.41%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "377.60"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.11%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0347"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "2.698.26"
$ws.Range("E45").Value = "  +0.58%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "133.48"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.93%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "23.81"
$c.ClearFormats()
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  -0.25%  "
